# Commit: Wed, Apr 01, 2020  5:07:27 PM
#
# The table on slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") has its
# table style switched from the deck's custom "Table_0" style
# ({F10F615E-BC91-4130-A38B-77058467BE7E}) to the built-in gallery style
# {FD6CA414-6D33-4445-A81E-44A27440F46F} (Table.Style itself is read-only
# in this object model, so ApplyStyle is used to make the write stick).

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{FD6CA414-6D33-4445-A81E-44A27440F46F}")
    }
}
